# Update the "想去人数" (F column) counts that changed when the gh-pages
# data was regenerated (commit: "Update gh-pages to output generated at 456a3b4").
#
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both contain the same
# rows of data (展览 = exhibitions, duplicated into the 全部类型 / all-types
# aggregate sheet), so the same underlying value updates show up twice.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 554
$ws1.Range("F6").Value  = 1146
$ws1.Range("F14").Value = 865
$ws1.Range("F22").Value = 2779
$ws1.Range("F25").Value = 2117
$ws1.Range("F27").Value = 2974
$ws1.Range("F28").Value = 566
$ws1.Range("F34").Value = 126

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 554
$ws4.Range("F6").Value  = 1146
$ws4.Range("F13").Value = 865
$ws4.Range("F22").Value = 2779
$ws4.Range("F27").Value = 2974
$ws4.Range("F28").Value = 566
$ws4.Range("F38").Value = 126
